$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new weekly price records (Perú Primera/Segunda avocado,
# $/bandeja 10 kilos) dated 2021-09-10 (serial 44449). They land right above
# the existing block of rows, so insert two blank rows at 237:238 first -
# this pushes the former rows 237-256 down to 239-258 intact.
$ws.Rows("237:238").Insert()

# Shared constant columns for this market/product across the whole block.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$tipo        = "Fruta"
$productoId  = 100106
$producto    = "Oleaginosos"
$categoriaId = 100106002
$categoria   = "Palta"
$variedad    = "Hass"
$unidad      = "`$/bandeja 10 kilos"
$origen      = "Perú"

# New row 237: Primera
$r = 237
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44449
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 360
$ws.Cells.Item($r, 14).Value = 25000
$ws.Cells.Item($r, 15).Value = 26000
$ws.Cells.Item($r, 16).Value = 25500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2550
$ws.Cells.Item($r, 20).Value = 10

# New row 238: Segunda
$r = 238
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44449
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 240
$ws.Cells.Item($r, 14).Value = 23000
$ws.Cells.Item($r, 15).Value = 24000
$ws.Cells.Item($r, 16).Value = 23500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2350
$ws.Cells.Item($r, 20).Value = 10
